# chart_of_accounts_mapping.xlsx - "tree structure for BS and IS"
# Fix a couple of typos and split the old "other assets" bucket into
# a "long term assets" subcategory (with "property and equipment" and
# "other assets" as its sub_category2 entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "cash and equivelents" -> "cash and equivalents"
$ws.Range("G2").Value2 = "cash and equivalents"

# Row 4 (account range 14700-16000): fix typo
# "property and equiptment" -> "property and equipment" first ...
$ws.Range("G4").Value2 = "property and equipment"
# ... then re-bucket the subcategory from "other assets" to "long term assets"
$ws.Range("F4").Value2 = "long term assets"

# Row 5 (account range 17500-19999): subcategory "other assets" becomes
# "long term assets"; add the now-distinct subcategory2 value "other assets"
$ws.Range("F5").Value2 = "long term assets"
$ws.Range("G5").Value2 = "other assets"

# Move the active selection to F6, matching the saved view state
$ws.Range("F6").Select()
